$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values for rows 2-5 (B and C columns switch from inline
# string "complex number" text to plain numbers; D4/D5 also change)
$ws.Range("B2").Value = -5
$ws.Range("C2").Value = 14.9990234375

$ws.Range("B3").Value = 40
$ws.Range("C3").Value = [double]"-1.20892581961463e+24"

$ws.Range("B4").Value = -5
$ws.Range("C4").Value = 14.9990234375
$ws.Range("D4").Value = 45

$ws.Range("B5").Value = -5
$ws.Range("C5").Value = 14.9990234375
$ws.Range("D5").Value = 0

# Remove the now-obsolete iteration rows 6-11
$ws.Range("A6:D11").ClearContents()
